$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "Intro, variables and arithmetic and assignment operators"
$ws.Range("E4").Value = "Comparison, identity and logical Operators"
$ws.Range("E5").Value = "Control flow (if, elif, else)"

$ws.Range("E4").Select()
